# Add the new Bacticure / Pulsera / Sognare Colchon Bioflex product rows
# to the "Hoja1" catalog sheet (rows 51-61, columns A:B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = "Familia de Producto" (product), Column B = "Product Category" (family)
$newRows = @(
    @("BACTICURE", "BACTICURE"),
    @("BACTICURE 1 FRASCO", "BACTICURE"),
    @("PULSERA", "PULSERA FORTUNARA"),
    @("PULSERA FORTUNARA", "PULSERA FORTUNARA"),
    @("SOGNARE COLCHON BIOFLEX KING", "SOGNARE COLCHON BIOFLEX"),
    @("SOGNARE COLCHON BIOFLEX INDIVIDUAL", "SOGNARE COLCHON BIOFLEX"),
    @("SOGNARE COLCHON BIOFLEX MATRIMONIAL", "SOGNARE COLCHON BIOFLEX"),
    @("SOGNARE COLCHON BIOFLEX MATRIMONIAL WEB", "SOGNARE COLCHON BIOFLEX"),
    @("SOGNARE COLCHON BIOFLEX QUEEN", "SOGNARE COLCHON BIOFLEX"),
    @("SOGNARE COLCHON BIOFLEX KING WEB", "SOGNARE COLCHON BIOFLEX"),
    @("SOGNARE COLCHON BIOFLEX QUEEN WEB", "SOGNARE COLCHON BIOFLEX")
)

$row = 51
foreach ($pair in $newRows) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Match the author's final view state: selection ends up on E49 in the
# saved workbook (near the newly appended rows).
$ws.Range("E49").Select()
